$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "target" column (G) values were shortened from "deuteron" to "d"
# for every data row (rows 2-11).
$ws.Range("G2:G11").Value = "d"

# Bold the header row and select it, matching the author's formatting pass.
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Bold = $true
$headerRange.Select()
